# Scheduled runner update: refresh market-price derived figures
# (currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns)
# across the per-job Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5273.25
$ws.Range("J70").Value = 6332.6665
$ws.Range("L70").Value = 18997.9995
$ws.Range("N70").Value = -19537.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5273.25
$ws.Range("J73").Value = 6332.6665
$ws.Range("L73").Value = 18997.9995
$ws.Range("N73").Value = -20869.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1896.1818
$ws.Range("I86").Value = 2383.8572
$ws.Range("K86").Value = 2383.8572
$ws.Range("M86").Value = -1260.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1896.1818
$ws.Range("I89").Value = 2383.8572
$ws.Range("K89").Value = 11919.286
$ws.Range("M89").Value = -6303.286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4586.1816
$ws.Range("I131").Value = 5162.222
$ws.Range("K131").Value = 15486.666
$ws.Range("M131").Value = -10446.666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1464.0212
$ws.Range("I132").Value = 1366.1818
$ws.Range("J132").Value = 2899
$ws.Range("K132").Value = 4098.5454
$ws.Range("L132").Value = 8697
$ws.Range("M132").Value = -1568.5454
$ws.Range("N132").Value = -13757

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2620.6707
$ws.Range("I138").Value = 922.8461
$ws.Range("K138").Value = 2768.5383
$ws.Range("M138").Value = 2371.4617

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1754.4
$ws.Range("I4").Value = 1363.6428
$ws.Range("J4").Value = 2666.1667
$ws.Range("K4").Value = 1363.6428
$ws.Range("L4").Value = 2666.1667
$ws.Range("M4").Value = -1247.6428
$ws.Range("N4").Value = -2898.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 465.75
$ws.Range("I5").Value = 326.27274
$ws.Range("K5").Value = 326.27274
$ws.Range("M5").Value = -214.27274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6267.4575
$ws.Range("I32").Value = 3885.1667
$ws.Range("K32").Value = 3885.1667
$ws.Range("M32").Value = -3598.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 52404.473
$ws.Range("I61").Value = 1699.8
$ws.Range("K61").Value = 1699.8
$ws.Range("M61").Value = -1487.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 53420.23
$ws.Range("I74").Value = 31773.666
$ws.Range("K74").Value = 31773.666
$ws.Range("M74").Value = -30899.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 53420.23
$ws.Range("I77").Value = 31773.666
$ws.Range("K77").Value = 158868.33
$ws.Range("M77").Value = -154500.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1384.7273
$ws.Range("I97").Value = 1155.7307
$ws.Range("K97").Value = 1155.7307
$ws.Range("M97").Value = -659.7307000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4126.826
$ws.Range("I102").Value = 3932.4856
$ws.Range("K102").Value = 3932.4856
$ws.Range("M102").Value = -2310.4856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 52404.473
$ws.Range("I136").Value = 1699.8
$ws.Range("K136").Value = 5099.4
$ws.Range("M136").Value = -2549.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 465.75
$ws.Range("I4").Value = 326.27274
$ws.Range("K4").Value = 326.27274
$ws.Range("M4").Value = -211.27274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2969.5356
$ws.Range("I20").Value = 2702.2273
$ws.Range("J20").Value = 3949.6667
$ws.Range("K20").Value = 2702.2273
$ws.Range("L20").Value = 3949.6667
$ws.Range("M20").Value = -2455.2273
$ws.Range("N20").Value = -4443.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 120792.07
$ws.Range("I86").Value = 75631.836
$ws.Range("J86").Value = 301433
$ws.Range("K86").Value = 75631.836
$ws.Range("L86").Value = 301433
$ws.Range("M86").Value = -74508.836
$ws.Range("N86").Value = -303679

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 120792.07
$ws.Range("I89").Value = 75631.836
$ws.Range("J89").Value = 301433
$ws.Range("K89").Value = 378159.18
$ws.Range("L89").Value = 1507165
$ws.Range("M89").Value = -372543.18
$ws.Range("N89").Value = -1518397

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1957.55
$ws.Range("I99").Value = 1635.375
$ws.Range("K99").Value = 1635.375
$ws.Range("M99").Value = -137.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 240.92308
$ws.Range("I7").Value = 85.42856999999999
$ws.Range("K7").Value = 85.42856999999999
$ws.Range("M7").Value = 27.57143000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6572.125
$ws.Range("I86").Value = 5172.643
$ws.Range("K86").Value = 5172.643
$ws.Range("M86").Value = -4049.643

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 6572.125
$ws.Range("I89").Value = 5172.643
$ws.Range("K89").Value = 25863.215
$ws.Range("M89").Value = -20247.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1491.4529
$ws.Range("I134").Value = 1346.4688
$ws.Range("K134").Value = 4039.4064
$ws.Range("M134").Value = -1504.4064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2399.5
$ws.Range("I113").Value = 873.6
$ws.Range("K113").Value = 2620.8
$ws.Range("M113").Value = -450.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.36
$ws.Range("I2").Value = 113.46154
$ws.Range("J2").Value = 148.66667
$ws.Range("K2").Value = 113.46154
$ws.Range("L2").Value = 148.66667
$ws.Range("M2").Value = -0.4615399999999994
$ws.Range("N2").Value = -374.66667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7286.2856
$ws.Range("I70").Value = 5252.5
$ws.Range("K70").Value = 5252.5
$ws.Range("M70").Value = -4982.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7286.2856
$ws.Range("I73").Value = 5252.5
$ws.Range("K73").Value = 5252.5
$ws.Range("M73").Value = -4316.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 18316.104
$ws.Range("I40").Value = 26346.834
$ws.Range("J40").Value = 5174.909
$ws.Range("K40").Value = 26346.834
$ws.Range("L40").Value = 5174.909
$ws.Range("M40").Value = -26210.834
$ws.Range("N40").Value = -5446.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3589.7
$ws.Range("I68").Value = 2362.375
$ws.Range("J68").Value = 8499
$ws.Range("K68").Value = 2362.375
$ws.Range("L68").Value = 8499
$ws.Range("M68").Value = -1613.375
$ws.Range("N68").Value = -9997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3589.7
$ws.Range("I71").Value = 2362.375
$ws.Range("J71").Value = 8499
$ws.Range("K71").Value = 11811.875
$ws.Range("L71").Value = 42495
$ws.Range("M71").Value = -8067.875
$ws.Range("N71").Value = -49983

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6656.684
$ws.Range("I132").Value = 5359.8887
$ws.Range("K132").Value = 16079.6661
$ws.Range("M132").Value = -13549.6661
